$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4171.493262381624
$ws.Range("C3").Value = 4066.41872927242
$ws.Range("C4").Value = 4039.985827795884
$ws.Range("C5").Value = 4039.985827795884
$ws.Range("C6").Value = 4039.985827795884
$ws.Range("C7").Value = 4039.985827795884
$ws.Range("C8").Value = 4039.985827795884
$ws.Range("C9").Value = 3957.067505419134
$ws.Range("C10").Value = 3957.067505419134
$ws.Range("C11").Value = 3951.996943790778
$ws.Range("C12").Value = 3951.996943790778
